$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-21 keep their existing timestamp (A) and label (B); only the
# sensor columns C:H are refreshed with the new window of readings.
$chData = @(
    @(4.275845370836725, -5.77583529645165, -10.19375841409568, -1.123454570770264, -0.7644603848457336, 2.25576400756836),
    @(0.6221833005047384, -4.507310633691367, -10.29080534301348, -0.0118510667234659, 1.250487327575684, 2.492785453796387),
    @(2.215511573241061, -1.993124782639061, -10.44278157637424, -1.767406940460205, 1.107741832733154, 2.145509243011475),
    @(4.329553655330059, 1.845118298626538, -8.935318364392899, -0.9109342694282532, 0.1809618026018142, 1.575592756271362),
    @(12.96200734976949, -5.231297192157526, 0.7199282710183832, -1.000949144363403, 0.0653806030750274, 1.502622127532959),
    @(-2.641056214403002, -7.670066910302078, -6.659932584570555, -1.359943270683289, 2.480335235595703, -0.1879525929689407),
    @(-10.01989353103133, -5.894216838299037, -4.208720184812602, 1.629188776016235, -4.939764499664307, -3.872169733047485),
    @(18.6824656044877, -5.135815050777897, -8.585992000246984, 4.826047420501709, 0.760465681552887, -4.821853160858154),
    @(20.39061579608261, -16.68983828461418, -10.62651905597464, 4.102732658386231, -0.6872287392616272, -7.78328800201416),
    @(37.93301056535422, -4.661358194863217, -0.6969916964537886, -5.43937349319458, -1.829724907875061, 0.9886320233345032),
    @(2.130085897125618, -8.977323400894459, -5.93280473491486, -8.410395622253418, 0.938364803791046, 10.93873500823975),
    @(-2.945472493267677, -38.83042448959086, -2.776809974004806, 3.432681560516357, 8.507601737976074, -2.705705165863037),
    @(4.761968791884873, -2.450371867058256, -21.97076380172833, 3.892343282699585, -3.907523155212401, 0.7548064589500427),
    @(-8.868669970723605, -16.28519311687272, -8.666557632036744, 10.68500232696533, -2.491254091262817, 0.1119193434715271),
    @(-7.487835359253294, -63.21479523581927, 15.10544187110538, -0.7697867155075073, 1.647830963134766, 1.185173273086548),
    @(6.800529531184445, -1.221064970797376, 2.408111162633729, -3.376488924026489, 1.38897168636322, 1.326853513717651),
    @(-10.6902578213072, -37.54249508748885, 29.03655345327873, -4.735766410827637, -0.1961418241262436, -5.656487464904785),
    @(-28.48631326624204, 16.53431056169849, -18.38939690429892, -2.693654298782349, 12.81020450592041, -6.878878593444824),
    @(-6.227460877207299, -5.338711252148522, -15.25327578487011, 4.974119186401367, -6.053498268127441, -4.147007942199707),
    @(-2.195728558021952, -9.250842510453845, -7.278287349931375, -2.26168966293335, -17.45715522766113, 0.2626542747020721)
)

$numRows = $chData.Count
$numCols = 6

$chArr = New-Object 'object[,]' $numRows, $numCols
for ($i = 0; $i -lt $numRows; $i++) {
    for ($j = 0; $j -lt $numCols; $j++) {
        $chArr[$i, $j] = $chData[$i][$j]
    }
}

$chRange = $ws.Range($ws.Cells.Item(2, 3), $ws.Cells.Item(2 + $numRows - 1, 8))
$chRange.Value = $chArr

# Rows 22-31 are brand new observations appended to the bottom of the
# table (timestamps 2000-2900), each with its own full A:H row.
$newData = @(
    @(2000, "walkingToRunning", 22.6993431820965, 9.916446260157828, 8.042439793580312, -1.576191902160644, 2.079262971878052, 1.673064351081848),
    @(2100, "walkingToRunning", 10.05702973692206, 4.430093637248776, 3.092578779130952, -10.3257417678833, 1.127981901168823, 5.69497013092041),
    @(2200, "walkingToRunning", -17.00681210844283, -24.96814476563599, -13.94656551924312, 6.999719619750977, -9.156081199645996, -6.504971027374268),
    @(2300, "walkingToRunning", 33.74133814101275, 5.435616121996308, -12.86083946612057, 1.670201539993286, -8.627177238464355, 0.0325571447610855),
    @(2400, "walkingToRunning", -5.288869998599139, -18.46977404460015, -9.24161248559119, 7.071092128753662, 1.176451444625854, -1.140298962593079),
    @(2500, "walkingToRunning", -25.29642678267203, -20.23905866738133, 10.65897996633644, -2.711231231689453, -4.331498622894287, 1.348158717155456),
    @(2600, "walkingToRunning", 3.936601254763914, 18.79409311281742, 10.11689840227166, -5.673731803894043, 1.553022384643555, 2.025134325027466),
    @(2700, "walkingToRunning", -36.19211913115269, -74.38502144013573, 47.36963078799701, -3.546931266784668, 9.443968772888184, -0.7339006066322327),
    @(2800, "walkingToRunning", -34.86098349654431, -1.257336098075001, -17.16005880880736, 0.6049370765686035, 8.170445442199707, -6.977948188781738),
    @(2900, "walkingToRunning", -23.05234499425689, -7.394026237846306, -17.08371260982228, 4.74082612991333, -5.404219627380371, -7.638944625854492)
)

$newNumRows = $newData.Count
$newNumCols = 8

$newArr = New-Object 'object[,]' $newNumRows, $newNumCols
for ($i = 0; $i -lt $newNumRows; $i++) {
    for ($j = 0; $j -lt $newNumCols; $j++) {
        $newArr[$i, $j] = $newData[$i][$j]
    }
}

$newStartRow = 22
$newEndRow = $newStartRow + $newNumRows - 1
$newRange = $ws.Range($ws.Cells.Item($newStartRow, 1), $ws.Cells.Item($newEndRow, 8))
$newRange.Value = $newArr
